$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A11").Value = "Tanita"
$ws.Range("B11").Value = "Menzel"
$ws.Range("C11").Value = "FEMALE"
$ws.Range("D11").Value = $false
